$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 406, shifting existing rows (406-477) down to (407-478)
$ws.Rows.Item(406).Insert()

# Populate the newly inserted row 406 with the new data record
$ws.Cells.Item(406, 1).Value = 1
$ws.Cells.Item(406, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(406, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(406, 4).Value = 45015
$ws.Cells.Item(406, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(406, 5).Value = 15
$ws.Cells.Item(406, 6).Value = 100114013
$ws.Cells.Item(406, 7).Value = "Zanahoria"
$ws.Cells.Item(406, 8).Value = "Sin especificar"
$ws.Cells.Item(406, 9).Value = "Primera"
$ws.Cells.Item(406, 10).Value = 25
$ws.Cells.Item(406, 11).Value = 24000
$ws.Cells.Item(406, 12).Value = 25000
$ws.Cells.Item(406, 13).Value = 24600
$ws.Cells.Item(406, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(406, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(406, 16).Value = 984
$ws.Cells.Item(406, 17).Value = 25
$ws.Cells.Item(406, 18).Value = "Hortaliza"

Write-Output "done"
